$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C to hold the PROJECT_ID data.
# This shifts the old C:F columns (and their formatting/widths) to D:G.
$ws.Columns("C:C").Insert()

# Copy the formatting of column B (style index used by every row, including
# the header) onto the freshly inserted column C so the new cells share the
# same font / alignment as the rest of the sheet.
$ws.Range("B1:B20").Copy()
$ws.Range("C1:C20").PasteSpecial(-4122)

# Give the new column a width close to the target (Excel quantizes column
# widths to the nearest 1/6 character, landing on 14.5 here).
$ws.Columns("C:C").ColumnWidth = 13.69

# Header for the new column.
$ws.Range("C1").Value = "PROJECT_ID"

# Fill in the example PROJECT_ID values. The values are entered grouped by
# project (Candy_kingdom, then Nightosphere, then Treehouse) so that the
# workbook's shared-string table records them in that same order.
$ws.Range("C2").Value = "Candy_kingdom"
$ws.Range("C3").Value = "Candy_kingdom"
$ws.Range("C4").Value = "Candy_kingdom"
$ws.Range("C5").Value = "Candy_kingdom"
$ws.Range("C6").Value = "Candy_kingdom"
$ws.Range("C7").Value = "Candy_kingdom"
$ws.Range("C8").Value = "Candy_kingdom"
$ws.Range("C15").Value = "Candy_kingdom"
$ws.Range("C16").Value = "Candy_kingdom"
$ws.Range("C17").Value = "Candy_kingdom"
$ws.Range("C18").Value = "Candy_kingdom"
$ws.Range("C19").Value = "Candy_kingdom"
$ws.Range("C20").Value = "Candy_kingdom"

$ws.Range("C11").Value = "Nightosphere"
$ws.Range("C12").Value = "Nightosphere"
$ws.Range("C14").Value = "Nightosphere"

$ws.Range("C9").Value = "Treehouse"
$ws.Range("C10").Value = "Treehouse"
$ws.Range("C13").Value = "Treehouse"

# Update the active cell/selection to match the saved view.
$ws.Range("F17").Select()
